# Auto-generated edit script: update crypto price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.759.53"
$ws.Range("E2").Value = "  -1.74%  "

$ws.Range("D3").Value = "1.942.79"
$ws.Range("E3").Value = "  -1.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.46%  "

$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4866"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.79%  "

$ws.Range("E8").Value = "  -2.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06852"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "105.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.88%  "

$ws.Range("D12").Value = "1.953.21"
$ws.Range("E12").Value = "  -0.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07746"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.304"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6951"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.16%  "

$ws.Range("D17").Value = "30.795.51"
$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007694"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.83%  "

$ws.Range("D20").Value = "2.208.27"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.433"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.02%  "

$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.447"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.685"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.161"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1037"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.387"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.552"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.524"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.350"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04833"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7453"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.64%  "

$ws.Range("E36").Value = "  -2.06%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01982"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.666"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.446"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "77.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.075"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.39%  "

$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4396"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9987"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.738"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.74%  "

$ws.Range("D48").Value = "1.001.40"
$ws.Range("E48").Value = "  +0.93%  "

$ws.Range("E49").Value = "  -2.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.115"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.18%  "
